$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in marks (value 5) for row 12: D12, E12, F12
$ws.Range("D12").Value = 5
$ws.Range("E12").Value = 5
$ws.Range("F12").Value = 5

# Fill in marks (value 5) for row 28: C28, D28, E28, F28
$ws.Range("C28").Value = 5
$ws.Range("D28").Value = 5
$ws.Range("E28").Value = 5
$ws.Range("F28").Value = 5

# Update frozen pane top-left cell and active selection to reflect scrolled view
$ws.Range("F12").Select()
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("C7").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("F12").Select()
